$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.970.36"
$ws.Range("E2").Value = "  -5.70%  "

$ws.Range("D3").Value = "'1.820.32"
$ws.Range("E3").Value = "  -5.40%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.80%  "

$ws.Range("D5").Value = "'327.37"
$ws.Range("E5").Value = "  -3.56%  "

$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").Value = "'0.4615"
$ws.Range("E7").Value = "  -4.18%  "

$ws.Range("D8").Value = "'0.3836"
$ws.Range("E8").Value = "  -5.38%  "

$ws.Range("D9").Value = "'45.88"
$ws.Range("E9").Value = "  -3.98%  "

$ws.Range("D10").Value = "'0.07802"
$ws.Range("E10").Value = "  -3.68%  "

$ws.Range("D11").Value = "'0.9562"

$ws.Range("D12").Value = "'21.82"
$ws.Range("E12").Value = "  -7.10%  "

$ws.Range("D13").Value = "'5.640"
$ws.Range("E13").Value = "  -5.90%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.765.14"
$ws.Range("E14").Value = "  -9.51%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.841"
$ws.Range("E15").Value = "  -4.90%  "

$ws.Range("D16").Value = "'0.06870"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").Value = "'86.31"
$ws.Range("E18").Value = "  -4.22%  "

$ws.Range("D19").Value = "'0.000009912"
$ws.Range("E19").Value = "  -3.60%  "

$ws.Range("D20").Value = "'16.74"
$ws.Range("E20").Value = "  -4.65%  "

$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -0.84%  "

$ws.Range("D22").Value = "'28.015.32"
$ws.Range("E22").Value = "  -5.59%  "

$ws.Range("D23").Value = "'5.308"

$ws.Range("D24").Value = "'10.91"
$ws.Range("E24").Value = "  -7.35%  "

$ws.Range("E25").Value = "  -1.90%  "

$ws.Range("D26").Value = "'1.970.54"
$ws.Range("E26").Value = "  -10.45%  "

$ws.Range("D27").Value = "'151.68"
$ws.Range("E27").Value = "  -3.29%  "

$ws.Range("D28").Value = "'19.13"
$ws.Range("E28").Value = "  -3.63%  "

$ws.Range("D29").Value = "'5.683"
$ws.Range("E29").Value = "  -13.60%  "

$ws.Range("D30").Value = "'1.963"
$ws.Range("E30").Value = "  -5.08%  "

$ws.Range("D31").Value = "'116.20"
$ws.Range("E31").Value = "  -3.44%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.09260"
$ws.Range("E32").Value = "  -3.46%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.9334"
$ws.Range("E33").Value = "  -6.75%  "

$ws.Range("D34").Value = "'5.269"
$ws.Range("E34").Value = "  -4.75%  "

$ws.Range("D35").Value = "'3.425"
$ws.Range("E35").Value = "  -3.33%  "

$ws.Range("D36").Value = "'1.303"
$ws.Range("E36").Value = "  -6.84%  "

$ws.Range("D37").Value = "'0.05942"
$ws.Range("E37").Value = "  -9.02%  "

$ws.Range("D38").Value = "'0.02142"
$ws.Range("E38").Value = "  -5.13%  "

$ws.Range("E39").Value = "  -4.33%  "

$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("D41").Value = "'7.528"
$ws.Range("E41").Value = "  -4.27%  "

$ws.Range("D42").Value = "'0.5572"
$ws.Range("E42").Value = "  -5.23%  "

$ws.Range("D43").Value = "'9.882"
$ws.Range("E43").Value = "  -7.26%  "

$ws.Range("D44").Value = "'0.1765"
$ws.Range("E44").Value = "  -3.61%  "

$ws.Range("D45").Value = "'1.219"
$ws.Range("E45").Value = "  -4.65%  "

$ws.Range("D46").Value = "'2.227"
$ws.Range("E46").Value = "  -9.72%  "

$ws.Range("D47").Value = "'11.55"
$ws.Range("E47").Value = "  -5.93%  "

$ws.Range("D48").Value = "'0.5236"
$ws.Range("E48").Value = "  -4.96%  "

$ws.Range("D49").Value = "'0.06993"
$ws.Range("E49").Value = "  -6.31%  "

$ws.Range("D50").Value = "'1.818"
$ws.Range("E50").Value = "  -7.35%  "

$ws.Range("D51").Value = "'112.24"
$ws.Range("E51").Value = "  -3.51%  "
